# Apply text-valued cell changes (company names, URLs, labels) for B, C, E columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textChanges = @(
    @{ Cell = 'B7'; Value = 'MXToken' },
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'E7'; Value = '6MXTokenMX' },
    @{ Cell = 'B8'; Value = 'FTXToken' },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' },
    @{ Cell = 'E8'; Value = '7FTXTokenFTT' },
    @{ Cell = 'B9'; Value = 'One' },
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' },
    @{ Cell = 'E9'; Value = '8OneONEBestin24h' },
    @{ Cell = 'B10'; Value = 'WazirX' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'E10'; Value = '9WazirXWRX' },
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'E11'; Value = '10MandalaExchangeTokenMDX' },
    @{ Cell = 'B12'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'E12'; Value = '11LiechtensteinCryptoassetsExchangeLCX' },
    @{ Cell = 'B13'; Value = 'BitrueCoin' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'E13'; Value = '12BitrueCoinBTR' },
    @{ Cell = 'B14'; Value = 'BitMartToken' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'E14'; Value = '13BitMartTokenBMX' },
    @{ Cell = 'B15'; Value = 'BitForexToken' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'E15'; Value = '14BitForexTokenBF' },
    @{ Cell = 'B16'; Value = 'CoinExToken' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' },
    @{ Cell = 'E16'; Value = '15CoinExTokenCET' },
    @{ Cell = 'B17'; Value = 'TigerCash' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'E17'; Value = '16TigerCashTCH' },
    @{ Cell = 'B18'; Value = 'HotbitToken' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb' },
    @{ Cell = 'E18'; Value = '17HotbitTokenHTB' },
    @{ Cell = 'B19'; Value = 'BitKan' },
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan' },
    @{ Cell = 'E19'; Value = '18BitKanKAN' },
    @{ Cell = 'B20'; Value = 'NitroEx' },
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx' },
    @{ Cell = 'E20'; Value = '19NitroExNTX' },
    @{ Cell = 'B21'; Value = 'LEO' },
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = 'E21'; Value = '20LEOLEO' },
    @{ Cell = 'B22'; Value = 'KuCoinToken' },
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs' },
    @{ Cell = 'E22'; Value = '21KuCoinTokenKCS' },
    @{ Cell = 'E41'; Value = '40KickTokenKICK' }
)

foreach ($change in $textChanges) {
    $ws.Range($change.Cell).Value = $change.Value
}

# Apply numeric-looking text values (Price column D) - must stay as TEXT, not be
# auto-converted to a number by Excel's type inference, and must not leave behind
# a lingering custom number-format style once done.
$priceChanges = @(
    @{ Cell = 'D2'; Value = '244.33' },
    @{ Cell = 'D3'; Value = '22.39' },
    @{ Cell = 'D4'; Value = '5.432' },
    @{ Cell = 'D5'; Value = '0.05766' },
    @{ Cell = 'D6'; Value = '3.427' },
    @{ Cell = 'D7'; Value = '0.8139' },
    @{ Cell = 'D8'; Value = '0.8821' },
    @{ Cell = 'D9'; Value = '0.01108' },
    @{ Cell = 'D10'; Value = '0.1441' },
    @{ Cell = 'D11'; Value = '0.07326' },
    @{ Cell = 'D12'; Value = '0.03019' },
    @{ Cell = 'D13'; Value = '0.03102' },
    @{ Cell = 'D14'; Value = '0.09406' },
    @{ Cell = 'D15'; Value = '0.001584' },
    @{ Cell = 'D16'; Value = '0.04828' },
    @{ Cell = 'D17'; Value = '0.006383' },
    @{ Cell = 'D18'; Value = '0.004138' },
    @{ Cell = 'D19'; Value = '0.0009946' },
    @{ Cell = 'D20'; Value = '0.0001500' },
    @{ Cell = 'D21'; Value = '3.725' },
    @{ Cell = 'D22'; Value = '6.320' },
    @{ Cell = 'D24'; Value = '0.3276' },
    @{ Cell = 'D25'; Value = '0.1320' },
    @{ Cell = 'D26'; Value = '4.179' },
    @{ Cell = 'D40'; Value = '0.03888' },
    @{ Cell = 'D41'; Value = '0.006789' },
    @{ Cell = 'D42'; Value = '0.1070' },
    @{ Cell = 'D43'; Value = '0.002801' },
    @{ Cell = 'D44'; Value = '0.006998' },
    @{ Cell = 'D45'; Value = '0.00005597' },
    @{ Cell = 'D48'; Value = '0.1559' }
)

foreach ($change in $priceChanges) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = '@'
    $cell.Value = $change.Value
    $cell.Style = 'Normal'
}
